$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos list refresh (GitHub Actions run).
# D-column prices are set as text (force "@" format then restore the
# default style) so numeric-looking strings such as "141.00" or
# "59.221.97" keep their exact original formatting instead of being
# auto-converted to floating point numbers by the COM Value setter.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.221.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.596.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.54%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.94%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.619.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.25%  "

$ws.Range("E11").Value = "  +2.21%  "

$ws.Range("E12").Value = "  +3.20%  "

$ws.Range("E13").Value = "  +2.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.055.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.115.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.635.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.26%  "

$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.406"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0726"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.38%  "

$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.44%  "

$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.80%  "

$ws.Range("E38").Value = "  +3.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.834"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "277.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.13%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0956"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.591"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.59%  "

$ws.Range("E47").Value = "  +1.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("E49").Value = "  +1.33%  "

$ws.Range("E50").Value = "  +2.53%  "

$ws.Range("E51").Value = "  +0.83%  "
